$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row at position 8 (new sample date 2022-10-13)
$ws.Rows(8).Insert()
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = 44847
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 100112040
$ws.Range("G8").Value = "Cilantro"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 7000
$ws.Range("N8").Value = "`$/caja 36 atados"
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 194
$ws.Range("Q8").Value = 36
$ws.Range("R8").Value = "Hortaliza"

# Insert new row at position 30 (after first insert shifts things down; new sample date 2022-10-14)
$ws.Rows(30).Insert()
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44848
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = 100112040
$ws.Range("G30").Value = "Cilantro"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 150
$ws.Range("K30").Value = 8000
$ws.Range("L30").Value = 8000
$ws.Range("M30").Value = 8000
$ws.Range("N30").Value = "`$/caja 36 atados"
$ws.Range("O30").Value = "Región del Maule"
$ws.Range("P30").Value = 222
$ws.Range("Q30").Value = 36
$ws.Range("R30").Value = "Hortaliza"

Write-Output "done"